$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.135.03"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.637.21"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "216.93"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").Value = "0.517"
$ws.Range("E6").Value = "  +1.43%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").Value = "20.10"
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "1.867.37"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "1.632.05"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "4.14"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "66.12"
$ws.Range("E16").Value = "  -2.67%  "
$ws.Range("D17").Value = "27.146.83"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "0.0₃0736"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "216.51"
$ws.Range("E19").Value = "  -2.46%  "
$ws.Range("D21").Value = "6.87"
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").Value = "147.58"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").Value = "15.67"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("D30").Value = "0.0508"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Value = "1.307.03"
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("D38").Value = "0.546"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").Value = "0.851"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("E41").Value = "  +4.22%  "
$ws.Range("D42").Value = "0.809"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").Value = "5.35"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").Value = "1.777.21"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").Value = "62.30"
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("D46").Value = "90.79"
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "7.63"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("D51").Value = "0.759"
$ws.Range("E51").Value = "  +13.86%  "
